$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New interviewers added to the Listing/Community/Household team (team 5),
# per corrections found from ToT review comments.
# Shared strings must be introduced column-by-column (all of column A, then
# column E) so the new unique strings land in the same order as the source
# workbook: shamimu/eline/joshua/lightness (names) before the Interview 0x
# comment strings.
$ws.Range("A20").Value = "shamimu"
$ws.Range("A21").Value = "eline"
$ws.Range("A22").Value = "joshua"
$ws.Range("A23").Value = "lightness"

$ws.Range("B20").Value = 5
$ws.Range("B21").Value = 5
$ws.Range("B22").Value = 5
$ws.Range("B23").Value = 5

$ws.Range("C20").Value = 800
$ws.Range("C21").Value = 900
$ws.Range("C22").Value = 850
$ws.Range("C23").Value = 950

$ws.Range("D20").Value = "interviewer"
$ws.Range("D21").Value = "interviewer"
$ws.Range("D22").Value = "interviewer"
$ws.Range("D23").Value = "interviewer"

$ws.Range("E20").Value = "Interview 01"
$ws.Range("E21").Value = "Interview 02"
$ws.Range("E22").Value = "Interview 03"
$ws.Range("E23").Value = "Interview 04"

# Leave the selection on the last-touched cell, matching the reviewer's
# final cursor position in the saved workbook.
$ws.Range("H21").Select()
